$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column (Price) cells keep their text formatting so Excel
# does not auto-convert numeric-looking strings (e.g. "1.00", "0.0320")
# into actual numbers and strip the trailing zeros / formatting.
$priceCells = @('D2', 'D3', 'D4', 'D5', 'D6', 'D7', 'D10', 'D13', 'D14', 'D15', 'D17', 'D19', 'D21', 'D22', 'D24', 'D25', 'D26', 'D27', 'D29', 'D30', 'D32', 'D33', 'D34', 'D35', 'D36', 'D40', 'D43', 'D44', 'D45', 'D46', 'D49', 'D50')
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated cell values from the source diff.
$ws.Range('D2').Value = '51.039.13'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '2.955.99'
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '381.00'
$ws.Range('E5').Value = '  +1.34%  '
$ws.Range('D6').Value = '102.11'
$ws.Range('E6').Value = '  -0.40%  '
$ws.Range('D7').Value = '0.546'
$ws.Range('E7').Value = '  +1.84%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +0.68%  '
$ws.Range('D10').Value = '36.52'
$ws.Range('E10').Value = '  -0.72%  '
$ws.Range('E11').Value = '  -0.76%  '
$ws.Range('E12').Value = '  +1.30%  '
$ws.Range('D13').Value = '12.47'
$ws.Range('E13').Value = '  +75.09%  '
$ws.Range('D14').Value = '18.40'
$ws.Range('E14').Value = '  +2.49%  '
$ws.Range('D15').Value = '3.416.63'
$ws.Range('E15').Value = '  +0.41%  '
$ws.Range('E16').Value = '  +5.29%  '
$ws.Range('D17').Value = '2.958.87'
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('E18').Value = '  +3.44%  '
$ws.Range('D19').Value = '51.096.08'
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('E20').Value = '  -2.45%  '
$ws.Range('D21').Value = '12.38'
$ws.Range('E21').Value = '  -1.59%  '
$ws.Range('D22').Value = '0.0₃0958'
$ws.Range('E22').Value = '  +0.21%  '
$ws.Range('E23').Value = '  +16.14%  '
$ws.Range('D24').Value = '269.37'
$ws.Range('E24').Value = '  +2.42%  '
$ws.Range('D25').Value = '69.71'
$ws.Range('D26').Value = '7.93'
$ws.Range('E26').Value = '  -2.34%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').Value = '  -0.97%  '
$ws.Range('D29').Value = '25.88'
$ws.Range('E29').Value = '  +0.81%  '
$ws.Range('D30').Value = '7.05'
$ws.Range('E30').Value = '  -11.06%  '
$ws.Range('E31').Value = '  -3.82%  '
$ws.Range('D32').Value = '10.42'
$ws.Range('E32').Value = '  +5.85%  '
$ws.Range('B33').Value = 'Toncoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D33').Value = '2.13'
$ws.Range('E33').Value = '  +5.40%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').Value = '51.12'
$ws.Range('E34').Value = '  +0.70%  '
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').Value = '34.24'
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('D36').Value = '0.0435'
$ws.Range('E36').Value = '  -4.91%  '
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('E38').Value = '  +9.01%  '
$ws.Range('E39').Value = '  +2.11%  '
$ws.Range('D40').Value = '16.71'
$ws.Range('E40').Value = '  +1.61%  '
$ws.Range('E41').Value = '  +2.68%  '
$ws.Range('E42').Value = '  -3.11%  '
$ws.Range('D43').Value = '124.55'
$ws.Range('E43').Value = '  +2.28%  '
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').Value = '3.56'
$ws.Range('E44').Value = '  +10.47%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '21.54'
$ws.Range('E45').Value = '  +2.11%  '
$ws.Range('D46').Value = '2.069.59'
$ws.Range('E46').Value = '  +3.46%  '
$ws.Range('E47').Value = '  -1.70%  '
$ws.Range('E48').Value = '  +1.39%  '
$ws.Range('D49').Value = '0.267'
$ws.Range('E49').Value = '  -2.39%  '
$ws.Range('D50').Value = '0.0320'
$ws.Range('E50').Value = '  -7.59%  '
$ws.Range('E51').Value = '  +6.47%  '

# Restore the default "Normal" style on the price cells so no extra
# style index is introduced (matches the original workbook which had
# no explicit style on these data cells).
foreach ($c in $priceCells) {
    $ws.Range($c).Style = "Normal"
}
